$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fullwidth/curly double quotes used by the author in the new line_id conditions.
$lq = [char]0x201C
$rq = [char]0x201D

# Update condition strings in column A (rows 5-7) to the corrected / quoted forms.
$ws.Cells.Item(5, 1).Value2 = 'status=1;point_type="8"'
$ws.Cells.Item(6, 1).Value2 = "status=1;line_id=$lq" + "1000000000056042" + "$rq"
$ws.Cells.Item(7, 1).Value2 = "status=1;line_id=$lq" + "1000000000056043" + "$rq"

# Row 7 (A7) gets a distinct font (等线) to set it apart, matching the new style
# introduced for the duplicated "普通地铁站深圳市龙华线" line.
$ws.Cells.Item(7, 1).Font.Name = "等线"

# The author's last selection moved from A4 to A7.
$ws.Range("A7").Select() | Out-Null
